$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "29÷7=4, 1"
$t.Cell(1, 2).Range.Text  = "31÷7=4, 3"
$t.Cell(1, 3).Range.Text  = "52÷2=26, 0"
$t.Cell(1, 4).Range.Text  = "19÷7=2, 5"
$t.Cell(1, 5).Range.Text  = "13÷4=3, 1"

$t.Cell(5, 1).Range.Text  = "29÷3=9, 2"
$t.Cell(5, 2).Range.Text  = "63÷8=7, 7"
$t.Cell(5, 3).Range.Text  = "93÷8=11, 5"
$t.Cell(5, 4).Range.Text  = "60÷8=7, 4"
$t.Cell(5, 5).Range.Text  = "33÷5=6, 3"

$t.Cell(9, 1).Range.Text  = "62÷8=7, 6"
$t.Cell(9, 2).Range.Text  = "32÷6=5, 2"
$t.Cell(9, 3).Range.Text  = "23÷3=7, 2"
$t.Cell(9, 4).Range.Text  = "49÷3=16, 1"
$t.Cell(9, 5).Range.Text  = "57÷9=6, 3"

$t.Cell(13, 1).Range.Text = "93÷9=10, 3"
$t.Cell(13, 2).Range.Text = "33÷7=4, 5"
$t.Cell(13, 3).Range.Text = "45÷4=11, 1"
$t.Cell(13, 4).Range.Text = "18÷7=2, 4"
$t.Cell(13, 5).Range.Text = "11÷9=1, 2"

$t.Cell(17, 1).Range.Text = "50÷7=7, 1"
$t.Cell(17, 2).Range.Text = "44÷3=14, 2"
$t.Cell(17, 3).Range.Text = "23÷8=2, 7"
$t.Cell(17, 4).Range.Text = "92÷6=15, 2"
$t.Cell(17, 5).Range.Text = "77÷2=38, 1"
